# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, reflecting a refreshed
# data pull (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (1st tab) ---
$ws1 = $wb.Worksheets.Item("展览")

$updates1 = @{
    2  = 583
    3  = 5538
    5  = 460
    9  = 1354
    12 = 3083
    13 = 1915
    14 = 119
    15 = 59
    17 = 34
    18 = 145
    21 = 352
    22 = 50
    23 = 3562
    24 = 1124
    25 = 2822
    26 = 283
    27 = 1988
    28 = 4087
    32 = 1304
    33 = 68
    35 = 1003
    36 = 1273
    38 = 1053
    39 = 681
    40 = 546
    41 = 414
    42 = 41
    44 = 3573
}

foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# --- Sheet "全部类型" (4th tab) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$updates4 = @{
    2  = 583
    3  = 5538
    9  = 1354
    10 = 3083
    12 = 1915
    13 = 119
    14 = 59
    18 = 34
    19 = 145
    21 = 352
    22 = 3562
    25 = 1124
    27 = 2822
    28 = 1988
    29 = 4087
    33 = 1304
    35 = 1003
    36 = 1273
    38 = 1053
    40 = 681
    42 = 414
    45 = 41
    48 = 3573
}

foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
